$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function CopyFormat($srcAddr, $dstAddr) {
  $ws.Range($srcAddr).Copy()
  $ws.Range($dstAddr).PasteSpecial(-4122)  # xlPasteFormats
}

$headers = @("IdCiclo","IdTorre","FechaInicio","FechaFin","TipoFin","CantidadNivelesCorrectos","PesoTorreFila","PesoTorreProducto","Lote","TiempoTotal")

# ----------------------------------------------------------------------
# Producto: Tarta de mousse de fresa (filas 103-108)
# ----------------------------------------------------------------------
CopyFormat "A92:B92" "A104:B104"
CopyFormat "A94" "A106"
CopyFormat "A95:J95" "A107:J107"
CopyFormat "C96:D96" "C108:D108"

$ws.Range("A104").Value = "Nombre: "
$ws.Range("B104").Value = "Tarta de mousse de fresa"
$ws.Range("A106").Value = "LISTA DE CICLOS"

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(107, $i + 1).Value = $headers[$i]
}

$ws.Range("A108").Value = 18
$ws.Range("B108").Value = 7
$ws.Range("C108").Value = 45830
$ws.Range("D108").Value = 45839.99998842592
$ws.Range("E108").Value = "Urgente"
$ws.Range("F108").Value = 3
$ws.Range("G108").Value = 250
$ws.Range("H108").Value = 750
$ws.Range("I108").Value = "L018"
$ws.Range("J108").Value = 6

$lo18 = $ws.ListObjects.Add(1, $ws.Range("A107:J108"), $null, 1)
$lo18.Name = "TablaCiclos_Tarta_de_mousse_de_fresa"

# ----------------------------------------------------------------------
# Producto: Panecillos de avena (filas 109-114)
# ----------------------------------------------------------------------
CopyFormat "A92:B92" "A110:B110"
CopyFormat "A94" "A112"
CopyFormat "A95:J95" "A113:J113"
CopyFormat "C96:D96" "C114:D114"

$ws.Range("A110").Value = "Nombre: "
$ws.Range("B110").Value = "Panecillos de avena"
$ws.Range("A112").Value = "LISTA DE CICLOS"

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(113, $i + 1).Value = $headers[$i]
}

$ws.Range("A114").Value = 19
$ws.Range("B114").Value = 8
$ws.Range("C114").Value = 45840
$ws.Range("D114").Value = 45849.99998842592
$ws.Range("E114").Value = "Normal"
$ws.Range("F114").Value = 2
$ws.Range("G114").Value = 70
$ws.Range("H114").Value = 140
$ws.Range("I114").Value = "L019"
$ws.Range("J114").Value = 1

$lo19 = $ws.ListObjects.Add(1, $ws.Range("A113:J114"), $null, 1)
$lo19.Name = "TablaCiclos_Panecillos_de_avena"

# ----------------------------------------------------------------------
# Ajustes de formato generales
# ----------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 25.16666666666667
